# Insert a new data row at row 503, shifting existing rows 503:557 down to 504:558,
# and populate the new row with the new data point.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(503).Insert()

$ws.Cells.Item(503, 1).Value = 11
$ws.Cells.Item(503, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(503, 3).Value = 'Bíobío'
$ws.Cells.Item(503, 4).Value = 45212
$ws.Cells.Item(503, 5).Value = 8
$ws.Cells.Item(503, 6).Value = 'Fruta'
$ws.Cells.Item(503, 7).Value = 100101
$ws.Cells.Item(503, 8).Value = 'Berries'
$ws.Cells.Item(503, 9).Value = 100112025
$ws.Cells.Item(503, 10).Value = 'Frutilla'
$ws.Cells.Item(503, 11).Value = 'Sin especificar'
$ws.Cells.Item(503, 12).Value = 'Especial'
$ws.Cells.Item(503, 13).Value = 100
$ws.Cells.Item(503, 14).Value = 11000
$ws.Cells.Item(503, 15).Value = 11000
$ws.Cells.Item(503, 16).Value = 11000
$ws.Cells.Item(503, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(503, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(503, 19).Value = 1571
$ws.Cells.Item(503, 20).Value = 7
